# Fruta / hortaliza, semanal
# Insert a new weekly price-sampling row at row 122 of the "Perejil" sheet
# (pushing the existing rows 122-154 down to 123-155) and populate it with
# the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 122, shifting rows 122:154 down
# to 123:155 (this also extends the used range to A1:R155, matching the
# updated <dimension> in the target file).
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new record's values.
$ws.Range("A122").Value = 8
$ws.Range("B122").Value = "Terminal La Palmera de La Serena"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 44754
$ws.Range("E122").Value = 4
$ws.Range("F122").Value = 100112044
$ws.Range("G122").Value = "Perejil"
$ws.Range("H122").Value = "Sin especificar"
$ws.Range("I122").Value = "Primera"
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 1500
$ws.Range("L122").Value = 2000
$ws.Range("M122").Value = 1750
$ws.Range("N122").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O122").Value = "Provincia del Elquí"
$ws.Range("P122").Value = 1167
$ws.Range("Q122").Value = 1.5
$ws.Range("R122").Value = "Hortaliza"
